$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.395.36'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.96%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.653.10'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -3.56%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.76'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9990'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.14%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -4.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '47.38'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -4.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3271'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -5.92%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.122'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -5.90%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06935'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -7.38%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.09%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.941'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -5.40%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.27'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -7.96%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.616'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -5.16%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.650.32'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -3.60%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001040'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -7.64%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06512'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -3.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9985'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.22%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '76.30'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -9.40%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.915'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -7.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.68'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -9.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.60'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.82%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.387.85'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.428'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.26%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.338'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -16.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '146.32'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -3.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.32'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -10.41%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.837.82'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.45%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.13'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -5.87%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.049'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.96%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.621'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -17.73%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08358'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -5.15%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.674'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -5.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.27'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -11.23%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.217'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -7.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06031'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -7.95%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02198'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -8.57%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -5.71%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2049'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -6.96%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.150'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -9.36%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9997'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5831'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -9.46%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.730'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.53%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.55'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -10.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5560'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -9.33%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '121.72'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -6.36%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.936'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -9.40%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06894'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -5.04%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '74.18'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -7.00%  '
